$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = '@'
$ws.Cells.Item(2,4).Value = '30.213.31'
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).NumberFormat = '@'
$ws.Cells.Item(2,5).Value = '  -0.63%  '
$ws.Cells.Item(2,5).Style = "Normal"

# Row 3
$ws.Cells.Item(3,4).NumberFormat = '@'
$ws.Cells.Item(3,4).Value = '1.882.03'
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).NumberFormat = '@'
$ws.Cells.Item(3,5).Value = '  -1.40%  '
$ws.Cells.Item(3,5).Style = "Normal"

# Row 4
$ws.Cells.Item(4,4).NumberFormat = '@'
$ws.Cells.Item(4,4).Value = '1.001'
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).NumberFormat = '@'
$ws.Cells.Item(4,5).Value = '  +0.16%  '
$ws.Cells.Item(4,5).Style = "Normal"

# Row 5
$ws.Cells.Item(5,4).NumberFormat = '@'
$ws.Cells.Item(5,4).Value = '237.37'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).NumberFormat = '@'
$ws.Cells.Item(5,5).Value = '  -0.46%  '
$ws.Cells.Item(5,5).Style = "Normal"

# Row 6
$ws.Cells.Item(6,4).NumberFormat = '@'
$ws.Cells.Item(6,4).Value = '1.003'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).NumberFormat = '@'
$ws.Cells.Item(6,5).Value = '  +0.31%  '
$ws.Cells.Item(6,5).Style = "Normal"

# Row 7
$ws.Cells.Item(7,4).NumberFormat = '@'
$ws.Cells.Item(7,4).Value = '0.4659'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).NumberFormat = '@'
$ws.Cells.Item(7,5).Value = '  -2.11%  '
$ws.Cells.Item(7,5).Style = "Normal"

# Row 8
$ws.Cells.Item(8,4).NumberFormat = '@'
$ws.Cells.Item(8,4).Value = '0.2799'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).NumberFormat = '@'
$ws.Cells.Item(8,5).Value = '  -2.09%  '
$ws.Cells.Item(8,5).Style = "Normal"

# Row 9
$ws.Cells.Item(9,4).NumberFormat = '@'
$ws.Cells.Item(9,4).Value = '0.06560'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).NumberFormat = '@'
$ws.Cells.Item(9,5).Value = '  -1.98%  '
$ws.Cells.Item(9,5).Style = "Normal"

# Row 10
$ws.Cells.Item(10,4).NumberFormat = '@'
$ws.Cells.Item(10,4).Value = '19.23'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).NumberFormat = '@'
$ws.Cells.Item(10,5).Value = '  +2.32%  '
$ws.Cells.Item(10,5).Style = "Normal"

# Row 11
$ws.Cells.Item(11,2).NumberFormat = '@'
$ws.Cells.Item(11,2).Value = 'Litecoin'
$ws.Cells.Item(11,2).Style = "Normal"
$ws.Cells.Item(11,3).NumberFormat = '@'
$ws.Cells.Item(11,3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(11,3).Style = "Normal"
$ws.Cells.Item(11,4).NumberFormat = '@'
$ws.Cells.Item(11,4).Value = '98.20'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).NumberFormat = '@'
$ws.Cells.Item(11,5).Value = '  -4.34%  '
$ws.Cells.Item(11,5).Style = "Normal"

# Row 12
$ws.Cells.Item(12,2).NumberFormat = '@'
$ws.Cells.Item(12,2).Value = 'TRON'
$ws.Cells.Item(12,2).Style = "Normal"
$ws.Cells.Item(12,3).NumberFormat = '@'
$ws.Cells.Item(12,3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12,3).Style = "Normal"
$ws.Cells.Item(12,4).NumberFormat = '@'
$ws.Cells.Item(12,4).Value = '0.07738'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).NumberFormat = '@'
$ws.Cells.Item(12,5).Value = '  +0.32%  '
$ws.Cells.Item(12,5).Style = "Normal"

# Row 13
$ws.Cells.Item(13,4).NumberFormat = '@'
$ws.Cells.Item(13,4).Value = '1.896.28'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).NumberFormat = '@'
$ws.Cells.Item(13,5).Value = '  -0.66%  '
$ws.Cells.Item(13,5).Style = "Normal"

# Row 14
$ws.Cells.Item(14,4).NumberFormat = '@'
$ws.Cells.Item(14,4).Value = '5.110'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).NumberFormat = '@'
$ws.Cells.Item(14,5).Value = '  -1.72%  '
$ws.Cells.Item(14,5).Style = "Normal"

# Row 15
$ws.Cells.Item(15,4).NumberFormat = '@'
$ws.Cells.Item(15,4).Value = '0.6593'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).NumberFormat = '@'
$ws.Cells.Item(15,5).Value = '  -2.40%  '
$ws.Cells.Item(15,5).Style = "Normal"

# Row 16
$ws.Cells.Item(16,4).NumberFormat = '@'
$ws.Cells.Item(16,4).Value = '283.88'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).NumberFormat = '@'
$ws.Cells.Item(16,5).Value = '  +9.58%  '
$ws.Cells.Item(16,5).Style = "Normal"

# Row 17
$ws.Cells.Item(17,4).NumberFormat = '@'
$ws.Cells.Item(17,4).Value = '30.190.15'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).NumberFormat = '@'
$ws.Cells.Item(17,5).Value = '  -0.74%  '
$ws.Cells.Item(17,5).Style = "Normal"

# Row 18
$ws.Cells.Item(18,4).NumberFormat = '@'
$ws.Cells.Item(18,4).Value = '1.002'
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).NumberFormat = '@'
$ws.Cells.Item(18,5).Value = '  +0.25%  '
$ws.Cells.Item(18,5).Style = "Normal"

# Row 19
$ws.Cells.Item(19,4).NumberFormat = '@'
$ws.Cells.Item(19,4).Value = '2.152.78'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).NumberFormat = '@'
$ws.Cells.Item(19,5).Value = '  -0.17%  '
$ws.Cells.Item(19,5).Style = "Normal"

# Row 20
$ws.Cells.Item(20,2).NumberFormat = '@'
$ws.Cells.Item(20,2).Value = 'Avalanche'
$ws.Cells.Item(20,2).Style = "Normal"
$ws.Cells.Item(20,3).NumberFormat = '@'
$ws.Cells.Item(20,3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(20,3).Style = "Normal"
$ws.Cells.Item(20,4).NumberFormat = '@'
$ws.Cells.Item(20,4).Value = '12.44'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).NumberFormat = '@'
$ws.Cells.Item(20,5).Value = '  -2.28%  '
$ws.Cells.Item(20,5).Style = "Normal"

# Row 21
$ws.Cells.Item(21,2).NumberFormat = '@'
$ws.Cells.Item(21,2).Value = 'ShibaInu'
$ws.Cells.Item(21,2).Style = "Normal"
$ws.Cells.Item(21,3).NumberFormat = '@'
$ws.Cells.Item(21,3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(21,3).Style = "Normal"
$ws.Cells.Item(21,4).NumberFormat = '@'
$ws.Cells.Item(21,4).Value = '0.000007278'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).NumberFormat = '@'
$ws.Cells.Item(21,5).Value = '  -3.00%  '
$ws.Cells.Item(21,5).Style = "Normal"

# Row 22
$ws.Cells.Item(22,4).NumberFormat = '@'
$ws.Cells.Item(22,4).Value = '5.311'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).NumberFormat = '@'
$ws.Cells.Item(22,5).Value = '  -1.84%  '
$ws.Cells.Item(22,5).Style = "Normal"

# Row 23
$ws.Cells.Item(23,2).NumberFormat = '@'
$ws.Cells.Item(23,2).Value = 'BinanceUSD'
$ws.Cells.Item(23,2).Style = "Normal"
$ws.Cells.Item(23,3).NumberFormat = '@'
$ws.Cells.Item(23,3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(23,3).Style = "Normal"
$ws.Cells.Item(23,4).NumberFormat = '@'
$ws.Cells.Item(23,4).Value = '1.000'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).NumberFormat = '@'
$ws.Cells.Item(23,5).Value = '  +0.08%  '
$ws.Cells.Item(23,5).Style = "Normal"

# Row 24
$ws.Cells.Item(24,4).NumberFormat = '@'
$ws.Cells.Item(24,4).Value = '6.174'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).NumberFormat = '@'
$ws.Cells.Item(24,5).Value = '  -1.96%  '
$ws.Cells.Item(24,5).Style = "Normal"

# Row 25
$ws.Cells.Item(25,2).NumberFormat = '@'
$ws.Cells.Item(25,2).Value = 'Cosmos'
$ws.Cells.Item(25,2).Style = "Normal"
$ws.Cells.Item(25,3).NumberFormat = '@'
$ws.Cells.Item(25,3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(25,3).Style = "Normal"
$ws.Cells.Item(25,4).NumberFormat = '@'
$ws.Cells.Item(25,4).Value = '9.218'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).NumberFormat = '@'
$ws.Cells.Item(25,5).Value = '  -2.49%  '
$ws.Cells.Item(25,5).Style = "Normal"

# Row 26
$ws.Cells.Item(26,2).NumberFormat = '@'
$ws.Cells.Item(26,2).Value = 'Monero'
$ws.Cells.Item(26,2).Style = "Normal"
$ws.Cells.Item(26,3).NumberFormat = '@'
$ws.Cells.Item(26,3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(26,3).Style = "Normal"
$ws.Cells.Item(26,4).NumberFormat = '@'
$ws.Cells.Item(26,4).Value = '165.67'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).NumberFormat = '@'
$ws.Cells.Item(26,5).Value = '  +0.77%  '
$ws.Cells.Item(26,5).Style = "Normal"

# Row 27
$ws.Cells.Item(27,4).NumberFormat = '@'
$ws.Cells.Item(27,4).Value = '18.94'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).NumberFormat = '@'
$ws.Cells.Item(27,5).Value = '  -0.13%  '
$ws.Cells.Item(27,5).Style = "Normal"

# Row 28
$ws.Cells.Item(28,4).NumberFormat = '@'
$ws.Cells.Item(28,4).Value = '1.997'
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).NumberFormat = '@'
$ws.Cells.Item(28,5).Value = '  -2.90%  '
$ws.Cells.Item(28,5).Style = "Normal"

# Row 29
$ws.Cells.Item(29,4).NumberFormat = '@'
$ws.Cells.Item(29,4).Value = '1.382'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).NumberFormat = '@'
$ws.Cells.Item(29,5).Value = '  +0.33%  '
$ws.Cells.Item(29,5).Style = "Normal"

# Row 30
$ws.Cells.Item(30,4).NumberFormat = '@'
$ws.Cells.Item(30,4).Value = '0.09813'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).NumberFormat = '@'
$ws.Cells.Item(30,5).Value = '  -2.76%  '
$ws.Cells.Item(30,5).Style = "Normal"

# Row 31
$ws.Cells.Item(31,4).NumberFormat = '@'
$ws.Cells.Item(31,4).Value = '4.463'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).NumberFormat = '@'
$ws.Cells.Item(31,5).Value = '  -3.20%  '
$ws.Cells.Item(31,5).Style = "Normal"

# Row 32
$ws.Cells.Item(32,4).NumberFormat = '@'
$ws.Cells.Item(32,4).Value = '1.493'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).NumberFormat = '@'
$ws.Cells.Item(32,5).Value = '  -1.33%  '
$ws.Cells.Item(32,5).Style = "Normal"

# Row 33
$ws.Cells.Item(33,4).NumberFormat = '@'
$ws.Cells.Item(33,4).Value = '4.178'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).NumberFormat = '@'
$ws.Cells.Item(33,5).Value = '  -1.86%  '
$ws.Cells.Item(33,5).Style = "Normal"

# Row 34
$ws.Cells.Item(34,4).NumberFormat = '@'
$ws.Cells.Item(34,4).Value = '0.04661'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).NumberFormat = '@'
$ws.Cells.Item(34,5).Value = '  -2.60%  '
$ws.Cells.Item(34,5).Style = "Normal"

# Row 35
$ws.Cells.Item(35,4).NumberFormat = '@'
$ws.Cells.Item(35,4).Value = '0.7068'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).NumberFormat = '@'
$ws.Cells.Item(35,5).Value = '  -3.32%  '
$ws.Cells.Item(35,5).Style = "Normal"

# Row 36
$ws.Cells.Item(36,4).NumberFormat = '@'
$ws.Cells.Item(36,4).Value = '1.089'
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).NumberFormat = '@'
$ws.Cells.Item(36,5).Value = '  -2.18%  '
$ws.Cells.Item(36,5).Style = "Normal"

# Row 37
$ws.Cells.Item(37,2).NumberFormat = '@'
$ws.Cells.Item(37,2).Value = 'Frax'
$ws.Cells.Item(37,2).Style = "Normal"
$ws.Cells.Item(37,3).NumberFormat = '@'
$ws.Cells.Item(37,3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(37,3).Style = "Normal"
$ws.Cells.Item(37,4).NumberFormat = '@'
$ws.Cells.Item(37,4).Value = '1.001'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).NumberFormat = '@'
$ws.Cells.Item(37,5).Value = '  +0.29%  '
$ws.Cells.Item(37,5).Style = "Normal"

# Row 38
$ws.Cells.Item(38,2).NumberFormat = '@'
$ws.Cells.Item(38,2).Value = 'HuobiToken'
$ws.Cells.Item(38,2).Style = "Normal"
$ws.Cells.Item(38,3).NumberFormat = '@'
$ws.Cells.Item(38,3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(38,3).Style = "Normal"
$ws.Cells.Item(38,4).NumberFormat = '@'
$ws.Cells.Item(38,4).Value = '2.717'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).NumberFormat = '@'
$ws.Cells.Item(38,5).Value = '  +0.32%  '
$ws.Cells.Item(38,5).Style = "Normal"

# Row 39
$ws.Cells.Item(39,2).NumberFormat = '@'
$ws.Cells.Item(39,2).Value = 'VeChain'
$ws.Cells.Item(39,2).Style = "Normal"
$ws.Cells.Item(39,3).NumberFormat = '@'
$ws.Cells.Item(39,3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39,3).Style = "Normal"
$ws.Cells.Item(39,4).NumberFormat = '@'
$ws.Cells.Item(39,4).Value = '0.01860'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).NumberFormat = '@'
$ws.Cells.Item(39,5).Value = '  -3.51%  '
$ws.Cells.Item(39,5).Style = "Normal"

# Row 40
$ws.Cells.Item(40,2).NumberFormat = '@'
$ws.Cells.Item(40,2).Value = 'FraxShare'
$ws.Cells.Item(40,2).Style = "Normal"
$ws.Cells.Item(40,3).NumberFormat = '@'
$ws.Cells.Item(40,3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(40,3).Style = "Normal"
$ws.Cells.Item(40,4).NumberFormat = '@'
$ws.Cells.Item(40,4).Value = '6.721'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).NumberFormat = '@'
$ws.Cells.Item(40,5).Value = '  +7.67%  '
$ws.Cells.Item(40,5).Style = "Normal"

# Row 41
$ws.Cells.Item(41,2).NumberFormat = '@'
$ws.Cells.Item(41,2).Value = 'MXToken'
$ws.Cells.Item(41,2).Style = "Normal"
$ws.Cells.Item(41,3).NumberFormat = '@'
$ws.Cells.Item(41,3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(41,3).Style = "Normal"
$ws.Cells.Item(41,4).NumberFormat = '@'
$ws.Cells.Item(41,4).Value = '2.522'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).NumberFormat = '@'
$ws.Cells.Item(41,5).Value = '  -3.01%  '
$ws.Cells.Item(41,5).Style = "Normal"

# Row 42
$ws.Cells.Item(42,2).NumberFormat = '@'
$ws.Cells.Item(42,2).Value = 'Aave'
$ws.Cells.Item(42,2).Style = "Normal"
$ws.Cells.Item(42,3).NumberFormat = '@'
$ws.Cells.Item(42,3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(42,3).Style = "Normal"
$ws.Cells.Item(42,4).NumberFormat = '@'
$ws.Cells.Item(42,4).Value = '72.46'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).NumberFormat = '@'
$ws.Cells.Item(42,5).Value = '  -3.18%  '
$ws.Cells.Item(42,5).Style = "Normal"

# Row 43
$ws.Cells.Item(43,2).NumberFormat = '@'
$ws.Cells.Item(43,2).Value = 'TrustWalletToken'
$ws.Cells.Item(43,2).Style = "Normal"
$ws.Cells.Item(43,3).NumberFormat = '@'
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(43,3).Style = "Normal"
$ws.Cells.Item(43,4).NumberFormat = '@'
$ws.Cells.Item(43,4).Value = '0.8691'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).NumberFormat = '@'
$ws.Cells.Item(43,5).Value = '  +0.39%  '
$ws.Cells.Item(43,5).Style = "Normal"

# Row 44
$ws.Cells.Item(44,2).NumberFormat = '@'
$ws.Cells.Item(44,2).Value = 'RenderToken'
$ws.Cells.Item(44,2).Style = "Normal"
$ws.Cells.Item(44,3).NumberFormat = '@'
$ws.Cells.Item(44,3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(44,3).Style = "Normal"
$ws.Cells.Item(44,4).NumberFormat = '@'
$ws.Cells.Item(44,4).Value = '1.944'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).NumberFormat = '@'
$ws.Cells.Item(44,5).Value = '  -2.34%  '
$ws.Cells.Item(44,5).Style = "Normal"

# Row 45
$ws.Cells.Item(45,2).NumberFormat = '@'
$ws.Cells.Item(45,2).Value = 'PaxDollar'
$ws.Cells.Item(45,2).Style = "Normal"
$ws.Cells.Item(45,3).NumberFormat = '@'
$ws.Cells.Item(45,3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(45,3).Style = "Normal"
$ws.Cells.Item(45,4).NumberFormat = '@'
$ws.Cells.Item(45,4).Value = '1.003'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).NumberFormat = '@'
$ws.Cells.Item(45,5).Value = '  +0.36%  '
$ws.Cells.Item(45,5).Style = "Normal"

# Row 46
$ws.Cells.Item(46,2).NumberFormat = '@'
$ws.Cells.Item(46,2).Value = 'Quant'
$ws.Cells.Item(46,2).Style = "Normal"
$ws.Cells.Item(46,3).NumberFormat = '@'
$ws.Cells.Item(46,3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(46,3).Style = "Normal"
$ws.Cells.Item(46,4).NumberFormat = '@'
$ws.Cells.Item(46,4).Value = '104.05'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).NumberFormat = '@'
$ws.Cells.Item(46,5).Value = '  -2.46%  '
$ws.Cells.Item(46,5).Style = "Normal"

# Row 47
$ws.Cells.Item(47,2).NumberFormat = '@'
$ws.Cells.Item(47,2).Value = 'TheSandbox'
$ws.Cells.Item(47,2).Style = "Normal"
$ws.Cells.Item(47,3).NumberFormat = '@'
$ws.Cells.Item(47,3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(47,3).Style = "Normal"
$ws.Cells.Item(47,4).NumberFormat = '@'
$ws.Cells.Item(47,4).Value = '0.4169'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).NumberFormat = '@'
$ws.Cells.Item(47,5).Value = '  -2.14%  '
$ws.Cells.Item(47,5).Style = "Normal"

# Row 48
$ws.Cells.Item(48,2).NumberFormat = '@'
$ws.Cells.Item(48,2).Value = 'Maker'
$ws.Cells.Item(48,2).Style = "Normal"
$ws.Cells.Item(48,3).NumberFormat = '@'
$ws.Cells.Item(48,3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(48,3).Style = "Normal"
$ws.Cells.Item(48,4).NumberFormat = '@'
$ws.Cells.Item(48,4).Value = '991.01'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).NumberFormat = '@'
$ws.Cells.Item(48,5).Value = '  -6.42%  '
$ws.Cells.Item(48,5).Style = "Normal"

# Row 49
$ws.Cells.Item(49,2).NumberFormat = '@'
$ws.Cells.Item(49,2).Value = 'Aptos'
$ws.Cells.Item(49,2).Style = "Normal"
$ws.Cells.Item(49,3).NumberFormat = '@'
$ws.Cells.Item(49,3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(49,3).Style = "Normal"
$ws.Cells.Item(49,4).NumberFormat = '@'
$ws.Cells.Item(49,4).Value = '7.200'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).NumberFormat = '@'
$ws.Cells.Item(49,5).Value = '  -3.78%  '
$ws.Cells.Item(49,5).Style = "Normal"

# Row 50
$ws.Cells.Item(50,2).NumberFormat = '@'
$ws.Cells.Item(50,2).Value = 'EnergySwap'
$ws.Cells.Item(50,2).Style = "Normal"
$ws.Cells.Item(50,3).NumberFormat = '@'
$ws.Cells.Item(50,3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50,3).Style = "Normal"
$ws.Cells.Item(50,4).NumberFormat = '@'
$ws.Cells.Item(50,4).Value = '9.099'
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).NumberFormat = '@'
$ws.Cells.Item(50,5).Value = '  +2.45%  '
$ws.Cells.Item(50,5).Style = "Normal"

# Row 51
$ws.Cells.Item(51,2).NumberFormat = '@'
$ws.Cells.Item(51,2).Value = 'Algorand'
$ws.Cells.Item(51,2).Style = "Normal"
$ws.Cells.Item(51,3).NumberFormat = '@'
$ws.Cells.Item(51,3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51,3).Style = "Normal"
$ws.Cells.Item(51,4).NumberFormat = '@'
$ws.Cells.Item(51,4).Value = '0.1164'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).NumberFormat = '@'
$ws.Cells.Item(51,5).Value = '  -3.01%  '
$ws.Cells.Item(51,5).Style = "Normal"
